$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Q0)
$ws.Range("B2").Value = 0.5692650327263563
$ws.Range("C2").Value = 0.5692650327263563
$ws.Range("D2").Value = 0.3908711553892542
$ws.Range("E2").Value = 0.6251968932978268
$ws.Range("F2").Value = 0.2682304822149452
$ws.Range("G2").Value = 14

# Row 3 (Q1)
$ws.Range("B3").Value = 0.3997456872501602
$ws.Range("C3").Value = 0.4033467742491634
$ws.Range("D3").Value = 0.2178465991588477
$ws.Range("E3").Value = 0.4667403980360472
$ws.Range("F3").Value = 0.2507737693766841
$ws.Range("G3").Value = 13

# Row 4 (Q2)
$ws.Range("B4").Value = 0.3400043916411057
$ws.Range("C4").Value = 0.3606492384221238
$ws.Range("D4").Value = 0.1784350047575419
$ws.Range("E4").Value = 0.4224156776891003
$ws.Range("F4").Value = 0.261809129132386
$ws.Range("G4").Value = 12

# Row 5 (Q3)
$ws.Range("B5").Value = 0.4288908803047028
$ws.Range("C5").Value = 0.4288908803047028
$ws.Range("D5").Value = 0.2376241711832545
$ws.Range("E5").Value = 0.4874670975391617
$ws.Range("F5").Value = 0.242990663137872
$ws.Range("G5").Value = 11

# Row 6 (Q4)
$ws.Range("B6").Value = 0.3767221034172891
$ws.Range("C6").Value = 0.3801506752190701
$ws.Range("D6").Value = 0.1848819966348901
$ws.Range("E6").Value = 0.429979065344919
$ws.Range("F6").Value = 0.2184858333361773
$ws.Range("G6").Value = 10

# Row 7 (Q5)
$ws.Range("B7").Value = 0.3516807317407905
$ws.Range("C7").Value = 0.358904738120446
$ws.Range("D7").Value = 0.1712159718543759
$ws.Range("E7").Value = 0.4137825175794355
$ws.Range("F7").Value = 0.2312546521125961
$ws.Range("G7").Value = 9

# Row 8 (Q6)
$ws.Range("B8").Value = 0.3826977698461033
$ws.Range("C8").Value = 0.3826977698461033
$ws.Range("D8").Value = 0.193400430583437
$ws.Range("E8").Value = 0.4397731580979414
$ws.Range("F8").Value = 0.2316224700134891
$ws.Range("G8").Value = 8

# Row 9 (Q7)
$ws.Range("B9").Value = 0.3590813353308283
$ws.Range("C9").Value = 0.3590813353308283
$ws.Range("D9").Value = 0.1682323819789578
$ws.Range("E9").Value = 0.4101614096657044
$ws.Range("F9").Value = 0.2141069966830871
$ws.Range("G9").Value = 7

# Row 10 (Q8)
$ws.Range("B10").Value = 0.3470367157308435
$ws.Range("C10").Value = 0.3470367157308435
$ws.Range("D10").Value = 0.1621328538698305
$ws.Range("E10").Value = 0.4026572411739673
$ws.Range("F10").Value = 0.2236918553848042
$ws.Range("G10").Value = 6

# Row 11 (Q9) - F11 is a new cell added
$ws.Range("B11").Value = 0.308333697360539
$ws.Range("C11").Value = 0.308333697360539
$ws.Range("D11").Value = 0.1232078055746762
$ws.Range("E11").Value = 0.3510096944169437
$ws.Range("F11").Value = 0.1875437837101504
$ws.Range("G11").Value = 5
